# Edit "type-of-data.xlsx": extend the controlled-vocabulary metadata block
# on sheet "Feuil2" with dct:creator / dct:publisher / owl:versionInfo /
# owl:versionIRI / owl:priorVersion / dct:license / bibo:status rows, update
# the dct:description text, and add the matching hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")

# Insert 7 new rows before the old header row (old row 4..10 were a blank
# gap followed by the "URI / skos:prefLabel / ..." table at row 7).  This
# pushes that table (and everything below it) down by 7 rows, so the old
# row 7 header becomes row 14, the three concept rows become 15-17 and the
# trailing formatted blank row becomes row 22 - matching the target layout.
$ws.Range("A4:A10").EntireRow.Insert()

# --- Row 2: dct:title (unchanged) ---

# --- Row 3: dct:description -- label unchanged, text updated further below ---

# --- Row 4: dct:creator ---
$ws.Range("A4").Value = "dct:creator"
$ws.Range("B4").Value = "Mario Scrocca (Cefriel)"

# --- Row 5: dct:publisher (value filled in further below) ---
$ws.Range("A5").Value = "dct:publisher"

# --- Row 6: owl:versionInfo ---
$ws.Range("A6").Value = "owl:versionInfo"
$ws.Range("B6").Value = "1.0.0"

# --- Row 7: owl:versionIRI (computed from the ConceptScheme URI + version) ---
$ws.Range("A7").Value = "owl:versionIRI"
$ws.Range("B7").Formula = "=CONCAT(B1,""/"",B6)"

# --- Row 8: owl:priorVersion label (left empty; hyperlink-style placeholder
#     is applied further below, after the "Hyperlink" style has already been
#     seeded via Hyperlinks.Add so it picks up the theme hyperlink color) ---
$ws.Range("A8").Value = "owl:priorVersion"

# --- Row 9: dct:license ---
$ws.Range("A9").Value = "dct:license"
$ws.Range("B9").Value = "https://creativecommons.org/licenses/by/4.0/"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://creativecommons.org/licenses/by/4.0/")

# Now apply the (seeded) "Hyperlink" style to the empty B8 placeholder cell.
$ws.Range("B8").Style = "Hyperlink"

# --- Row 10: bibo:status ---
$ws.Range("A10").Value = "http://purl.org/ontology/bibo/status"
$ws.Range("B10").Value = "Published Controlled Vocabulary"
$ws.Range("C10").Clear()

# --- Row 3: dct:description text (added after the rows below, matching the
#     original authoring order) ---
$ws.Range("B3").Value = "Controlled vocabulary for the type of data defined in TANGENT."

# --- Row 5: dct:publisher value (added last) ---
$ws.Range("B5").Value = "TANGENT WP2"

# --- Row 1: ConceptScheme URI hyperlink (added after the license hyperlink,
#     matching the original authoring order) ---
$ws.Hyperlinks.Add($ws.Range("B1"), "https://knowledge.c-innovationhub.com/tangent/type-of-data")

# Put the selection where the author left it.
$ws.Range("A21").Select()
